$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update depression response text to add suicide hotline note (rows 69-73)
$ws.Range("C69:C73").Value = "Depression, otherwise known as Major Depressive Disorder is an unfortunately common and serious mental illness that negatively affects how a person feels, thinks, and acts. Depression causes overwhelming feelings of sadness and/or a loss of interest in activities you may have once enjoyed. Depression also can lead to a wide variety of physical problems and can decrease your ability to function at work, school, and at home. Some symptoms of depression include feeling sad or having a depressed mood, loss of interest or pleasure in your hobbies, changes in appetite, and many more. For more information please visit ""https://www.psychiatry.org/patients-families/depression/what-is-depression"".  If you think are experiencing depression, please consult your primary care physician or psychiatrist. If you are having thoughts of harming yourself in any way, please text or call the National Suicide Hotline at (988)."

# Rows 99-103: replace "what is suicidal ideation" rows with bipolar disorder content
$ws.Range("A99").Value = "what is bipolar"
$ws.Range("A100").Value = "can you tell me what bipolar disorder is"
$ws.Range("A101").Value = "what are the signs of bipolar"
$ws.Range("A102").Value = "tell me some  of the symptoms of bipolar"
$ws.Range("A103").Value = "can you explain what bipolar disorder is"
$ws.Range("B99:B103").Value = "what is bipolar disorder"

# Rows 104-108: expand "how to apologize" prompts (category unchanged)
$ws.Range("A104").Value = "what is the best way to apologize to someone"
$ws.Range("A105").Value = "how can I tell them im sorry"
$ws.Range("A106").Value = "how do I apologize to them"
$ws.Range("A107").Value = "what can I do to apologize"
$ws.Range("A108").Value = "what is the best way to say im sorry"

# Rows 109-113: expand "how to build confidence" prompts (category moved from "how to calm down")
$ws.Range("A109").Value = "how can I work on being less nervous around people"
$ws.Range("A110").Value = "what can I do to be more confident"
$ws.Range("A111").Value = "how can I be more confident in myself"
$ws.Range("A112").Value = "what are some ways I can be more confident around others"
$ws.Range("A113").Value = "tell me some ways to improve my confidence"
$ws.Range("B109:B113").Value = "how to build confidence"

# Rows 114-118: "how to cheer someone else up" prompts moved up (category moved from "how to ground")
$ws.Range("A114").Value = "How can I cheer someone up?"
$ws.Range("A115").Value = "Can you help me cheer up someone?"
$ws.Range("A116").Value = "Tell me how to cheer someone up"
$ws.Range("A117").Value = "What are some ways to cheer someone up"
$ws.Range("A118").Value = "what ways can I cheer people up"
$ws.Range("B114:B118").Value = "how to cheer someone else up"

# Rows 119-123: "how is sai" prompts moved up (category moved from "how to cope")
$ws.Range("A119").Value = "How are you"
$ws.Range("A120").Value = "how are you feeling today"
$ws.Range("A121").Value = "tell me how you are feeling"
$ws.Range("A122").Value = "are you okay"
$ws.Range("A123").Value = "are you feeling okay"
$ws.Range("B119:B123").Value = "how is sai"
$ws.Range("C119:C123").Value = "I'm feeling great today! How are you feeling?"

# Rows 124-128: category changes to "what is a panic attack" (was "how to make friends")
$ws.Range("B124:B128").Value = "what is a panic attack"

# Rows 129-133: category changes to "how to calm down" (was "how to build confidence")
$ws.Range("B129:B133").Value = "how to calm down"

# Rows 134-138: category changes to "how to cope" (was "how to be less nervous")
$ws.Range("B134:B138").Value = "how to cope"

# Rows 139-143: clear old prompts, category changes to "how to make friends" (was "how to cheer someone else up")
$ws.Range("A139:A143").ClearContents()
$ws.Range("B139:B143").Value = "how to make friends"

# Rows 144-148 removed entirely (suicidewatch.csv source data dropped)
$ws.Range("A144:C148").ClearContents()

# Restore view state
$ws.Range("E71").Select()
